# Apply the refreshed cryptocurrency data (price + 1h volume change) that
# GitHub Actions pulled for cryptos.xlsx: updated Price/Volume(1h) values
# per row, plus a 3-way reshuffle of the RenderToken/MXToken/FraxShare rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) sometimes hold values that *look* numeric
# ("1.000", "0.4615", ...) but must be stored as literal text, exactly as
# authored (Excel would otherwise parse them as numbers and drop trailing
# zeros / normalize them). Force text format on those cells before writing,
# then restore the workbook's original (unstyled/General) formatting.
$priceTextCellRefs = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($ref in $priceTextCellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.756.50"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "1.799.18"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "308.92"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4615"
$ws.Range("E7").Value = "  +3.44%  "
$ws.Range("D8").Value = "0.3717"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "0.07265"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "0.8554"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("D11").Value = "20.40"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").Value = "1.763.94"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "5.316"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "6.494"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "0.07045"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "90.71"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "0.000008633"
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "14.63"
$ws.Range("E20").Value = "  -4.03%  "
$ws.Range("D21").Value = "26.747.69"
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").Value = "5.289"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").Value = "2.004.14"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "1.908"
$ws.Range("E25").Value = "  -4.47%  "
$ws.Range("D26").Value = "150.37"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "18.20"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "2.138"
$ws.Range("E28").Value = "  -10.07%  "
$ws.Range("D29").Value = "5.220"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").Value = "114.07"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("D31").Value = "0.08903"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").Value = "0.7559"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("D33").Value = "1.161"
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").Value = "4.432"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").Value = "2.889"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "1.121"
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("D38").Value = "0.01941"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "0.05207"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.899"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.173"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.361"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("D43").Value = "0.5219"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").Value = "0.1648"
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("D45").Value = "8.508"
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("D46").Value = "0.5003"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").Value = "10.22"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").Value = "104.15"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "1.648"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("E51").Value = "  -1.37%  "

# Restore default (General) styling on the Price cells we text-formatted.
foreach ($ref in $priceTextCellRefs) {
    $ws.Range($ref).Style = "Normal"
}

